# Add "Training Voucher" program as a new row in the programs worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 18

$ws.Cells.Item($row, 1).Value = "trainingVoucher"
$ws.Cells.Item($row, 2).Value = "Training Voucher"
$ws.Cells.Item($row, 3).Value = 2003
$ws.Cells.Item($row, 4).Value = "Active Labor Market Policy"
$ws.Cells.Item($row, 5).Value = 39.03
$ws.Cells.Item($row, 6).Value = "Training Vouchers ""Bildungsgutscheine"" are awarded to unemployed to allow them to participate in training programs usually lasting between several months and 3 years."
$ws.Cells.Item($row, 7).Value = 1901.21

$ws.Range("G19").Select()
